$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new rows at the top of the data block for the new listings
# added in this week's push (KB제27호스팩, 유진테크놀로지, 유투바이오,
# 퀄리타스반도체, 워트).
$ws.Rows("2:6").Insert()

# Force plain text on the date columns (A-C) and the rate column (O) so
# Excel doesn't auto-coerce date-like strings ("2023-10-19") into date
# serials, or percentage-like strings ("17.63%") into numeric percentages,
# before we write the literal text values below. The numeric columns are
# left alone so they stay real numbers.
$ws.Range("A2:C6").NumberFormat = "@"
$ws.Range("O2:O6").NumberFormat = "@"

$ws.Range("A2").Value = "2023-10-19"
$ws.Range("B2").Value = "2023-10-20"
$ws.Range("C2").Value = "2023-11-03"
$ws.Range("D2").Value = "KB"
$ws.Range("E2").Value = "KB제27호스팩"
$ws.Range("F2").Value = 12500000
$ws.Range("G2").Value = 12500000
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 2000
$ws.Range("J2").Value = 2000
$ws.Range("K2").Value = 12905000
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 2000
$ws.Range("N2").Value = "10.14:1"
$ws.Range("O2").Value = "-"
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = 0
$ws.Range("U2").Value = 0
$ws.Range("V2").Value = 0
$ws.Range("W2").Value = 0
$ws.Range("X2").Value = 0
$ws.Range("Y2").Value = "기업인수합병"

$ws.Range("A3").Value = "2023-10-11"
$ws.Range("B3").Value = "2023-10-17"
$ws.Range("C3").Value = "2023-11-02"
$ws.Range("D3").Value = "NH"
$ws.Range("E3").Value = "유진테크놀로지"
$ws.Range("F3").Value = 1049482
$ws.Range("G3").Value = 944534
$ws.Range("H3").Value = 104948
$ws.Range("I3").Value = 12800
$ws.Range("J3").Value = 14500
$ws.Range("K3").Value = 6261485
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 17000
$ws.Range("N3").Value = "914.02:1"
$ws.Range("O3").Value = "17.63%"
$ws.Range("P3").Value = 34557425427
$ws.Range("Q3").Value = 39824841246
$ws.Range("R3").Value = 23231897516
$ws.Range("S3").Value = 3011651602
$ws.Range("T3").Value = 2384643399
$ws.Range("U3").Value = 4436005255
$ws.Range("V3").Value = 2755379556
$ws.Range("W3").Value = 1176755354
$ws.Range("X3").Value = 4046949430
$ws.Range("Y3").Value = "이차전지 정밀금형 외"

$ws.Range("A4").Value = "2023-10-18"
$ws.Range("B4").Value = "2023-10-19"
$ws.Range("C4").Value = "2023-11-02"
$ws.Range("D4").Value = "신한"
$ws.Range("E4").Value = "유투바이오"
$ws.Range("F4").Value = 1128720
$ws.Range("G4").Value = 1128720
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 3300
$ws.Range("J4").Value = 3900
$ws.Range("K4").Value = 11287196
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 4400
$ws.Range("N4").Value = "1,276.73:1"
$ws.Range("O4").Value = "2.35%"
$ws.Range("P4").Value = 50552623684
$ws.Range("Q4").Value = 69013134090
$ws.Range("R4").Value = 16887814423
$ws.Range("S4").Value = 9027232647
$ws.Range("T4").Value = 15707796256
$ws.Range("U4").Value = 555562560
$ws.Range("V4").Value = 9236341465
$ws.Range("W4").Value = 13159994846
$ws.Range("X4").Value = 1236029732
$ws.Range("Y4").Value = "체외진단검사서비스, 의료IT솔루션"

$ws.Range("A5").Value = "2023-10-06"
$ws.Range("B5").Value = "2023-10-13"
$ws.Range("C5").Value = "2023-10-27"
$ws.Range("D5").Value = "한국"
$ws.Range("E5").Value = "퀄리타스반도체"
$ws.Range("F5").Value = 1800000
$ws.Range("G5").Value = 1800000
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 13000
$ws.Range("J5").Value = 15000
$ws.Range("K5").Value = 10193520
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 17000
$ws.Range("N5").Value = "818.52:1"
$ws.Range("O5").Value = "22.25%"
$ws.Range("P5").Value = 3952070068
$ws.Range("Q5").Value = 10789274729
$ws.Range("R5").Value = 6040367765
$ws.Range("S5").Value = -4741577598
$ws.Range("T5").Value = -3671026788
$ws.Range("U5").Value = -3292521363
$ws.Range("V5").Value = -3654889267
$ws.Range("W5").Value = -2281074471
$ws.Range("X5").Value = -2539505707
$ws.Range("Y5").Value = "초고속 통신용 반도체 IP"

$ws.Range("A6").Value = "2023-10-05"
$ws.Range("B6").Value = "2023-10-11"
$ws.Range("C6").Value = "2023-10-26"
$ws.Range("D6").Value = "키움"
$ws.Range("E6").Value = "워트"
$ws.Range("F6").Value = 4000000
$ws.Range("G6").Value = 4000000
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 5000
$ws.Range("J6").Value = 5600
$ws.Range("K6").Value = 16120000
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = 6500
$ws.Range("N6").Value = "793.26:1"
$ws.Range("O6").Value = "10.06%"
$ws.Range("P6").Value = 26691070764
$ws.Range("Q6").Value = 22835113396
$ws.Range("R6").Value = 7519073723
$ws.Range("S6").Value = 6041270235
$ws.Range("T6").Value = 6705120210
$ws.Range("U6").Value = 982044017
$ws.Range("V6").Value = 5119066139
$ws.Range("W6").Value = 5870093710
$ws.Range("X6").Value = 1198587470
$ws.Range("Y6").Value = "초정밀 온습도 제어장비"

# Strip the temporary formatting so the new rows end up with the same
# "no explicit style" look as the rest of the data rows.
$ws.Range("A2:Y6").ClearFormats()

# Drop the 4 oldest listings that fell off the bottom of the table
# (originally rows 13-16: 대신밸런스제16호스팩, 유안타제11호스팩,
# 한국제12호스팩, 대신밸런스제15호스팩) - after the insert above they
# now sit at rows 18-21.
$ws.Rows("18:21").Delete()
